$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '24.681.20'
$ws.Range("E2").Value = '  +2.49%  '

# Row 3
$ws.Range("D3").Value = '1.696.18'
$ws.Range("E3").Value = '  +1.54%  '

# Row 4
$ws.Range("D4").Value = '''1.000'
$ws.Range("E4").Value = '  +0.27%  '

# Row 5
$ws.Range("D5").Value = '''316.44'
$ws.Range("E5").Value = '  +2.11%  '

# Row 6
$ws.Range("D6").Value = '''0.9991'
$ws.Range("E6").Value = '  +0.13%  '

# Row 7
$ws.Range("D7").Value = '''0.3950'
$ws.Range("E7").Value = '  +1.77%  '

# Row 8
$ws.Range("D8").Value = '''0.4039'
$ws.Range("E8").Value = '  +1.09%  '

# Row 9
$ws.Range("D9").Value = '''1.516'
$ws.Range("E9").Value = '  +2.80%  '

# Row 10
$ws.Range("D10").Value = '''0.9987'
$ws.Range("E10").Value = '  +0.13%  '

# Row 11
$ws.Range("D11").Value = '''51.22'
$ws.Range("E11").Value = '  -4.25%  '

# Row 12
$ws.Range("D12").Value = '''0.08792'
$ws.Range("E12").Value = '  +1.26%  '

# Row 13
$ws.Range("D13").Value = '''7.222'
$ws.Range("E13").Value = '  +5.97%  '

# Row 14
$ws.Range("D14").Value = '''23.48'
$ws.Range("E14").Value = '  +3.17%  '

# Row 15
$ws.Range("D15").Value = '''8.204'
$ws.Range("E15").Value = '  +12.11%  '

# Row 16
$ws.Range("D16").Value = '''0.00001322'
$ws.Range("E16").Value = '  +0.93%  '

# Row 17
$ws.Range("D17").Value = '1.693.88'
$ws.Range("E17").Value = '  +1.58%  '

# Row 18
$ws.Range("D18").Value = '''99.95'
$ws.Range("E18").Value = '  +0.76%  '

# Row 19
$ws.Range("D19").Value = '''0.07007'
$ws.Range("E19").Value = '  +1.30%  '

# Row 20
$ws.Range("D20").Value = '''19.76'
$ws.Range("E20").Value = '  +3.12%  '

# Row 21
$ws.Range("D21").Value = '''7.089'
$ws.Range("E21").Value = '  +7.40%  '

# Row 22
$ws.Range("D22").Value = '''0.9993'
$ws.Range("E22").Value = '  +0.25%  '

# Row 23
$ws.Range("D23").Value = '''14.34'
$ws.Range("E23").Value = '  +3.57%  '

# Row 24
$ws.Range("D24").Value = '24.672.61'
$ws.Range("E24").Value = '  +2.45%  '

# Row 25
$ws.Range("D25").Value = '''3.136'
$ws.Range("E25").Value = '  +3.24%  '

# Row 26
$ws.Range("E26").Value = '  +1.60%  '

# Row 27
$ws.Range("D27").Value = '''22.84'
$ws.Range("E27").Value = '  +4.83%  '

# Row 28
$ws.Range("D28").Value = '''161.96'
$ws.Range("E28").Value = '  +0.95%  '

# Row 29
$ws.Range("D29").Value = '''137.57'
$ws.Range("E29").Value = '  +5.44%  '

# Row 30
$ws.Range("D30").Value = '''5.203'
$ws.Range("E30").Value = '  +1.28%  '

# Row 31
$ws.Range("D31").Value = '''7.456'
$ws.Range("E31").Value = '  +2.38%  '

# Row 32
$ws.Range("D32").Value = '1.879.46'

# Row 33
$ws.Range("E33").Value = '  -1.51%  '

# Row 34
$ws.Range("D34").Value = '''0.08628'
$ws.Range("E34").Value = '  +0.16%  '

# Row 35
$ws.Range("D35").Value = '''7.112'
$ws.Range("E35").Value = '  -0.85%  '

# Row 36
$ws.Range("D36").Value = '''11.57'
$ws.Range("E36").Value = '  +6.99%  '

# Row 37
$ws.Range("D37").Value = '''0.2755'
$ws.Range("E37").Value = '  +3.85%  '

# Row 38
$ws.Range("D38").Value = '''1.926'
$ws.Range("E38").Value = '  +0.32%  '

# Row 39
$ws.Range("D39").Value = '''14.50'
$ws.Range("E39").Value = '  -0.47%  '

# Row 40
$ws.Range("D40").Value = '''0.09242'
$ws.Range("E40").Value = '  +4.68%  '

# Row 41
$ws.Range("D41").Value = '''0.02727'
$ws.Range("E41").Value = '  +7.51%  '

# Row 42
$ws.Range("D42").Value = '''1.482'
$ws.Range("E42").Value = '  +2.53%  '

# Row 43
$ws.Range("D43").Value = '''0.7677'
$ws.Range("E43").Value = '  +1.49%  '

# Row 44
$ws.Range("D44").Value = '''2.656'
$ws.Range("E44").Value = '  +9.79%  '

# Row 45
$ws.Range("E45").Value = '  +5.19%  '

# Row 46
$ws.Range("D46").Value = '''0.7189'
$ws.Range("E46").Value = '  +1.34%  '

# Row 47
$ws.Range("E47").Value = '  +2.72%  '

# Row 48
$ws.Range("E48").Value = '  +0.19%  '

# Row 49
$ws.Range("B49").Value = 'Flow'
$ws.Range("C49").Value = 'https://coinranking.com/coin/QQ0NCmjVq+flow-flow'
$ws.Range("D49").Value = '''1.333'
$ws.Range("E49").Value = '  +8.62%  '

# Row 50
$ws.Range("B50").Value = 'Quant'
$ws.Range("C50").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D50").Value = '''141.00'
$ws.Range("E50").Value = '  +1.42%  '

# Row 51
$ws.Range("D51").Value = '''0.07993'
$ws.Range("E51").Value = '  +1.94%  '
